$d = $word.ActiveDocument

$replacements = @(
    @{old="926÷7=132, 2"; new="640÷6=106, 4"},
    @{old="486÷4=121, 2"; new="300÷2=150, 0"},
    @{old="873÷9=97, 0"; new="419÷6=69, 5"},
    @{old="696÷3=232, 0"; new="620÷8=77, 4"},
    @{old="269÷8=33, 5"; new="496÷6=82, 4"},
    @{old="179÷2=89, 1"; new="855÷7=122, 1"},
    @{old="831÷7=118, 5"; new="195÷6=32, 3"},
    @{old="723÷5=144, 3"; new="961÷3=320, 1"},
    @{old="360÷2=180, 0"; new="915÷9=101, 6"},
    @{old="987÷8=123, 3"; new="626÷5=125, 1"},
    @{old="539÷7=77, 0"; new="472÷9=52, 4"},
    @{old="778÷8=97, 2"; new="719÷6=119, 5"},
    @{old="118÷8=14, 6"; new="894÷9=99, 3"},
    @{old="619÷8=77, 3"; new="938÷3=312, 2"},
    @{old="746÷6=124, 2"; new="757÷4=189, 1"},
    @{old="397÷6=66, 1"; new="606÷6=101, 0"},
    @{old="350÷8=43, 6"; new="579÷5=115, 4"},
    @{old="222÷9=24, 6"; new="215÷7=30, 5"},
    @{old="632÷3=210, 2"; new="129÷3=43, 0"},
    @{old="884÷5=176, 4"; new="342÷2=171, 0"},
    @{old="852÷2=426, 0"; new="107÷7=15, 2"},
    @{old="650÷5=130, 0"; new="465÷6=77, 3"},
    @{old="224÷7=32, 0"; new="634÷4=158, 2"},
    @{old="980÷7=140, 0"; new="221÷7=31, 4"},
    @{old="265÷2=132, 1"; new="555÷4=138, 3"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "done"
